$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the anchor paragraph ("I've started to research Hibernate...")
# and from it work out the index of the blank paragraph that precedes
# the existing "_GoBack" bookmark paragraph. All the new content is
# spliced in right after that blank paragraph (i.e. right before the
# bookmark paragraph).
# ------------------------------------------------------------------
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "research Hibernate and how to implement") {
        $idx = $i
        break
    }
}

$cursor = $idx + 1   # index of the blank paragraph right after the anchor

# --- New paragraph: hyperlink to the Hibernate research paper ---
$d.Paragraphs.Item($cursor).Range.InsertParagraphAfter()
$cursor = $cursor + 1
$p1 = $d.Paragraphs.Item($cursor)
$hlRange = $p1.Range
$hlRange.Collapse(1)
$hibernateUrl = "https://github.com/Costa-Eurico/CIT-360-W16-Understanding-Portfolio/blob/master/Understanding%20Portfolio%20Submission%2002/Hibernate/Hibernate%20Research%20Paper.docx"
$null = $d.Hyperlinks.Add($hlRange, $hibernateUrl, "", "", $hibernateUrl)

# --- Empty paragraph ---
$d.Paragraphs.Item($cursor).Range.InsertParagraphAfter()
$cursor = $cursor + 1

# --- New paragraph: "The following code is my first implementation..." ---
$d.Paragraphs.Item($cursor).Range.InsertParagraphAfter()
$cursor = $cursor + 1
$p3 = $d.Paragraphs.Item($cursor)
$p3.Range.Text = "The following code is my first implementation of the hibernate framework and is tailored for my personal app project."

# --- Empty paragraph ---
$d.Paragraphs.Item($cursor).Range.InsertParagraphAfter()
$cursor = $cursor + 1

# --- New paragraph: two runs about past frameworks / ESB work ---
$d.Paragraphs.Item($cursor).Range.InsertParagraphAfter()
$cursor = $cursor + 1
$p5 = $d.Paragraphs.Item($cursor)
$part1 = "I have in the past worked on similar frameworks such as Hibernate. A few years ago, when I was working heavily with Microsoft .NET, and there were no similar frameworks around, I created my own framework, initially tailored for a project I was leading for a large Insurance Company back in Portugal. At the time, we wanted to create a web services layer in front of Siebel CRM 7.5, and basically creating a simple ESB to handle any type of transactions coming into Siebel, regardless of the format and transport used (SOAP web services, RESTful APIs, queues, etc.). the configuration for the transactions, as well as xslt transformations and maps between the incoming and outgoing transactions, and the canonical model we created were all in a MS SQL database, and instead of making this ESB tightly coupled to the database, I created a framework that abstracted the "
$part2 = "ESB from the database, so that we could easily make the ESB more portable and depending on the customer" + [char]0x2019 + "s database preference, use either Oracle, or MS SQL, or anything else. We ended up creating support for both MSSQL and Oracle DB at the time. "
$p5.Range.Text = $part1
$p5again = $d.Paragraphs.Item($cursor)
$p5again.Range.InsertAfter($part2)

# --- Empty paragraph ---
$d.Paragraphs.Item($cursor).Range.InsertParagraphAfter()
$cursor = $cursor + 1

# --- Final text goes INSIDE the existing bookmark paragraph, right
#     before the bookmark itself, so the bookmark ("_GoBack") stays
#     intact and attached to the end of this paragraph, exactly as in
#     the target document. ---
$bookmarkParaIdx = $cursor + 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIdx)
$bookmarkRange = $bookmarkPara.Range
$bookmarkRange.Collapse(1)
$closingText = "So, this to explain that I understand the importance of having frameworks such as Hibernate to abstract these types of details from the application" + [char]0x2019 + "s core. Developers spend way too much time writing code to do the same things such as database access, over and over. These frameworks remove the need to do so."
$bookmarkRange.InsertBefore($closingText)

Write-Host "Hibernate narrative inserted."
